$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H33").Value = 91.59999999999999
$ws.Range("I33").Value = 115
$ws.Range("J33").Value = 56.5
$ws.Range("K33").Value = 115
$ws.Range("L33").Value = 56.5
$ws.Range("M33").Value = 114
$ws.Range("N33").Value = -514.5
$ws.Range("H62").Value = 4625
$ws.Range("I62").Value = 4229.3335
$ws.Range("J62").Value = 4921.75
$ws.Range("K62").Value = 4229.3335
$ws.Range("L62").Value = 4921.75
$ws.Range("M62").Value = -3605.3335
$ws.Range("N62").Value = -6169.75
$ws.Range("H65").Value = 4625
$ws.Range("I65").Value = 4229.3335
$ws.Range("J65").Value = 4921.75
$ws.Range("K65").Value = 21146.6675
$ws.Range("L65").Value = 24608.75
$ws.Range("M65").Value = -18026.6675
$ws.Range("N65").Value = -30848.75
$ws.Range("H113").Value = 2091.7693
$ws.Range("J113").Value = 2091.7693
$ws.Range("L113").Value = 2091.7693
$ws.Range("N113").Value = -8599.7693
$ws.Range("H116").Value = 2961375.5
$ws.Range("I116").Value = 3665698.2
$ws.Range("J116").Value = 3219.4
$ws.Range("K116").Value = 3665698.2
$ws.Range("L116").Value = 3219.4
$ws.Range("M116").Value = -3662256.2
$ws.Range("N116").Value = -10103.4
$ws.Range("H134").Value = 41796
$ws.Range("J134").Value = 41796
$ws.Range("L134").Value = 41796
$ws.Range("N134").Value = -51936
$ws.Range("H137").Value = 93680.45
$ws.Range("I137").Value = 4497.3335
$ws.Range("J137").Value = 127124.125
$ws.Range("K137").Value = 13492.0005
$ws.Range("L137").Value = 381372.375
$ws.Range("M137").Value = -10942.0005
$ws.Range("N137").Value = -386472.375
$ws.Range("H138").Value = 3034.1636
$ws.Range("J138").Value = 4827.5
$ws.Range("L138").Value = 14482.5
$ws.Range("N138").Value = -24762.5

$ws = $wb.Worksheets.Item(2)
$ws.Range("H74").Value = 5024.625
$ws.Range("I74").Value = 5923.8423
$ws.Range("J74").Value = 1607.6
$ws.Range("K74").Value = 5923.8423
$ws.Range("L74").Value = 1607.6
$ws.Range("M74").Value = -5049.8423
$ws.Range("N74").Value = -3355.6
$ws.Range("H77").Value = 5024.625
$ws.Range("I77").Value = 5923.8423
$ws.Range("J77").Value = 1607.6
$ws.Range("K77").Value = 29619.2115
$ws.Range("L77").Value = 8038
$ws.Range("M77").Value = -25251.2115
$ws.Range("N77").Value = -16774
$ws.Range("H97").Value = 540.04
$ws.Range("I97").Value = 536.875
$ws.Range("J97").Value = 545.6667
$ws.Range("K97").Value = 536.875
$ws.Range("L97").Value = 545.6667
$ws.Range("M97").Value = -40.875
$ws.Range("N97").Value = -1537.6667
$ws.Range("H102").Value = 1844.1666
$ws.Range("I102").Value = 1442.5714
$ws.Range("J102").Value = 2406.4
$ws.Range("K102").Value = 1442.5714
$ws.Range("L102").Value = 2406.4
$ws.Range("M102").Value = 179.4286
$ws.Range("N102").Value = -5650.4
$ws.Range("H122").Value = 2263999
$ws.Range("I122").Value = 2942648.8
$ws.Range("J122").Value = 1833.3334
$ws.Range("K122").Value = 8827946.399999999
$ws.Range("L122").Value = 5500.0002
$ws.Range("M122").Value = -8825496.399999999
$ws.Range("N122").Value = -10400.0002

$ws = $wb.Worksheets.Item(3)
$ws.Range("H64").Value = 387.16666
$ws.Range("I64").Value = 507.875
$ws.Range("J64").Value = 290.6
$ws.Range("K64").Value = 507.875
$ws.Range("L64").Value = 290.6
$ws.Range("M64").Value = -282.875
$ws.Range("N64").Value = -740.6
$ws.Range("H67").Value = 387.16666
$ws.Range("I67").Value = 507.875
$ws.Range("J67").Value = 290.6
$ws.Range("K67").Value = 507.875
$ws.Range("L67").Value = 290.6
$ws.Range("M67").Value = 272.125
$ws.Range("N67").Value = -1850.6
$ws.Range("H80").Value = 95.90909000000001
$ws.Range("I80").Value = 132.4
$ws.Range("J80").Value = 65.5
$ws.Range("K80").Value = 132.4
$ws.Range("L80").Value = 65.5
$ws.Range("M80").Value = 865.6
$ws.Range("N80").Value = -2061.5
$ws.Range("H83").Value = 95.90909000000001
$ws.Range("I83").Value = 132.4
$ws.Range("J83").Value = 65.5
$ws.Range("K83").Value = 662
$ws.Range("L83").Value = 327.5
$ws.Range("M83").Value = 4330
$ws.Range("N83").Value = -10311.5

$ws = $wb.Worksheets.Item(4)
$ws.Range("H94").Value = 3350.4092
$ws.Range("I94").Value = 1885.4445
$ws.Range("J94").Value = 4364.615
$ws.Range("K94").Value = 1885.4445
$ws.Range("L94").Value = 4364.615
$ws.Range("M94").Value = -1434.4445
$ws.Range("N94").Value = -5266.615
$ws.Range("H99").Value = 2867.6365
$ws.Range("I99").Value = 2380.6875
$ws.Range("J99").Value = 4166.1665
$ws.Range("K99").Value = 2380.6875
$ws.Range("L99").Value = 4166.1665
$ws.Range("M99").Value = -882.6875
$ws.Range("N99").Value = -7162.1665
$ws.Range("H126").Value = 2867.6365
$ws.Range("I126").Value = 2380.6875
$ws.Range("J126").Value = 4166.1665
$ws.Range("K126").Value = 7142.0625
$ws.Range("L126").Value = 12498.4995
$ws.Range("M126").Value = -4672.0625
$ws.Range("N126").Value = -17438.4995

$ws = $wb.Worksheets.Item(5)
$ws.Range("H5").Value = 681
$ws.Range("I5").Value = 681
$ws.Range("K5").Value = 2043
$ws.Range("M5").Value = -1931
$ws.Range("H80").Value = 494.5
$ws.Range("J80").Value = 494.5
$ws.Range("L80").Value = 1483.5
$ws.Range("N80").Value = -3355.5
$ws.Range("H83").Value = 494.5
$ws.Range("J83").Value = 494.5
$ws.Range("L83").Value = 4450.5
$ws.Range("N83").Value = -13810.5
$ws.Range("H122").Value = 728.1579
$ws.Range("J122").Value = 1001.5
$ws.Range("L122").Value = 9013.5
$ws.Range("N122").Value = -13913.5
$ws.Range("H124").Value = 2959.2222
$ws.Range("I124").Value = 800
$ws.Range("J124").Value = 3576.1428
$ws.Range("K124").Value = 2400
$ws.Range("L124").Value = 10728.4284
$ws.Range("M124").Value = 2510
$ws.Range("N124").Value = -20548.4284
$ws.Range("H132").Value = 674131.9399999999
$ws.Range("I132").Value = 633.5454999999999
$ws.Range("J132").Value = 2526252.5
$ws.Range("K132").Value = 5701.9095
$ws.Range("L132").Value = 22736272.5
$ws.Range("M132").Value = -3171.9095
$ws.Range("N132").Value = -22741332.5
$ws.Range("H135").Value = 681
$ws.Range("I135").Value = 681
$ws.Range("K135").Value = 6129
$ws.Range("M135").Value = -3594

$ws = $wb.Worksheets.Item(6)
$ws.Range("H102").Value = 2012.5333
$ws.Range("I102").Value = 1908
$ws.Range("J102").Value = 2300
$ws.Range("K102").Value = 1908
$ws.Range("L102").Value = 2300
$ws.Range("M102").Value = -286
$ws.Range("N102").Value = -5544
$ws.Range("H107").Value = 472.41666
$ws.Range("I107").Value = 308.625
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 308.625
$ws.Range("L107").Value = 800
$ws.Range("M107").Value = 1611.375
$ws.Range("N107").Value = -4640
$ws.Range("H113").Value = 8458.25
$ws.Range("I113").Value = 1952.2858
$ws.Range("J113").Value = 54000
$ws.Range("K113").Value = 1952.2858
$ws.Range("L113").Value = 54000
$ws.Range("M113").Value = 217.7141999999999
$ws.Range("N113").Value = -58340
$ws.Range("H122").Value = 3450792.2
$ws.Range("I122").Value = 4763873.5
$ws.Range("J122").Value = 3954.125
$ws.Range("K122").Value = 14291620.5
$ws.Range("L122").Value = 11862.375
$ws.Range("M122").Value = -14289170.5
$ws.Range("N122").Value = -16762.375
$ws.Range("H132").Value = 3350.389
$ws.Range("I132").Value = 2820.7334
$ws.Range("J132").Value = 5998.6665
$ws.Range("K132").Value = 8462.200199999999
$ws.Range("L132").Value = 17995.9995
$ws.Range("M132").Value = -5932.200199999999
$ws.Range("N132").Value = -23055.9995

$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 1697.591
$ws.Range("I7").Value = 1437.7273
$ws.Range("J7").Value = 1957.4546
$ws.Range("K7").Value = 1437.7273
$ws.Range("L7").Value = 1957.4546
$ws.Range("M7").Value = -1325.7273
$ws.Range("N7").Value = -2181.4546
$ws.Range("H22").Value = 770223.3
$ws.Range("I22").Value = 1111544.5
$ws.Range("J22").Value = 2250.5
$ws.Range("K22").Value = 1111544.5
$ws.Range("L22").Value = 2250.5
$ws.Range("M22").Value = -1111249.5
$ws.Range("N22").Value = -2840.5
$ws.Range("H27").Value = 770223.3
$ws.Range("I27").Value = 1111544.5
$ws.Range("J27").Value = 2250.5
$ws.Range("K27").Value = 1111544.5
$ws.Range("L27").Value = 2250.5
$ws.Range("M27").Value = -1111437.5
$ws.Range("N27").Value = -2464.5
$ws.Range("H40").Value = 2959.926
$ws.Range("I40").Value = 2020.95
$ws.Range("J40").Value = 5642.7144
$ws.Range("K40").Value = 2020.95
$ws.Range("L40").Value = 5642.7144
$ws.Range("M40").Value = -1884.95
$ws.Range("N40").Value = -5914.7144
$ws.Range("H61").Value = 6381.85
$ws.Range("I61").Value = 7102.4707
$ws.Range("J61").Value = 2298.3333
$ws.Range("K61").Value = 7102.4707
$ws.Range("L61").Value = 2298.3333
$ws.Range("M61").Value = -6900.4707
$ws.Range("N61").Value = -2702.3333
$ws.Range("H113").Value = 6381.85
$ws.Range("I113").Value = 7102.4707
$ws.Range("J113").Value = 2298.3333
$ws.Range("K113").Value = 7102.4707
$ws.Range("L113").Value = 2298.3333
$ws.Range("M113").Value = -4932.4707
$ws.Range("N113").Value = -6638.3333
$ws.Range("H126").Value = 1697.591
$ws.Range("I126").Value = 1437.7273
$ws.Range("J126").Value = 1957.4546
$ws.Range("K126").Value = 4313.1819
$ws.Range("L126").Value = 5872.3638
$ws.Range("M126").Value = -1843.1819
$ws.Range("N126").Value = -10812.3638
$ws.Range("H132").Value = 10477.6
$ws.Range("I132").Value = 13463.333
$ws.Range("K132").Value = 40389.999
$ws.Range("M132").Value = -37859.999

$ws = $wb.Worksheets.Item(8)
$ws.Range("H122").Value = 3657.5
$ws.Range("I122").Value = 3100
$ws.Range("J122").Value = 3769
$ws.Range("K122").Value = 9300
$ws.Range("L122").Value = 11307
$ws.Range("M122").Value = -6850
$ws.Range("N122").Value = -16207
$ws.Range("H126").Value = 792.3333
$ws.Range("I126").Value = 855.0909
$ws.Range("J126").Value = 619.75
$ws.Range("K126").Value = 2565.2727
$ws.Range("L126").Value = 1859.25
$ws.Range("M126").Value = -95.27269999999999
$ws.Range("N126").Value = -6799.25
